$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1295.5714
$ws.Range("I18").Value = 594
$ws.Range("J18").Value = 3049.5
$ws.Range("K18").Value = 594
$ws.Range("L18").Value = 3049.5
$ws.Range("M18").Value = -310
$ws.Range("N18").Value = -3617.5

$ws.Range("H43").Value = 418001.4
$ws.Range("I43").Value = 20000
$ws.Range("J43").Value = 683335.7
$ws.Range("K43").Value = 20000
$ws.Range("L43").Value = 683335.7
$ws.Range("M43").Value = -19931
$ws.Range("N43").Value = -683473.7

$ws.Range("H68").Value = 85195
$ws.Range("I68").Value = 29990
$ws.Range("J68").Value = 103596.664
$ws.Range("K68").Value = 29990
$ws.Range("L68").Value = 103596.664
$ws.Range("M68").Value = -29241
$ws.Range("N68").Value = -105094.664

$ws.Range("H71").Value = 85195
$ws.Range("I71").Value = 29990
$ws.Range("J71").Value = 103596.664
$ws.Range("K71").Value = 89970
$ws.Range("L71").Value = 310789.992
$ws.Range("M71").Value = -86226
$ws.Range("N71").Value = -318277.992

$ws.Range("H86").Value = 1542409.2
$ws.Range("I86").Value = 2502915
$ws.Range("J86").Value = 5600
$ws.Range("K86").Value = 2502915
$ws.Range("L86").Value = 5600
$ws.Range("M86").Value = -2501792

$ws.Range("H89").Value = 1542409.2
$ws.Range("I89").Value = 2502915
$ws.Range("J89").Value = 5600
$ws.Range("K89").Value = 12514575
$ws.Range("L89").Value = 28000
$ws.Range("M89").Value = -12508959

$ws.Range("H112").Value = 1785.1
$ws.Range("I112").Value = 1046
$ws.Range("J112").Value = 2031.4667
$ws.Range("K112").Value = 3138
$ws.Range("L112").Value = 6094.4001
$ws.Range("M112").Value = -2030

$ws.Range("H137").Value = 5516.815
$ws.Range("I137").Value = 3910.5264
$ws.Range("J137").Value = 9331.75
$ws.Range("K137").Value = 11731.5792
$ws.Range("L137").Value = 27995.25
$ws.Range("M137").Value = -9181.5792

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 0
$ws.Range("N62").Value = $null  # remove (was -41248)

$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 0
$ws.Range("N65").Value = $null  # remove (was -126240)

$ws.Range("H74").Value = 7432431
$ws.Range("I74").Value = 10871875
$ws.Range("J74").Value = 840163
$ws.Range("K74").Value = 10871875
$ws.Range("L74").Value = 840163
$ws.Range("M74").Value = -10871001

$ws.Range("H77").Value = 7432431
$ws.Range("I77").Value = 10871875
$ws.Range("J77").Value = 840163
$ws.Range("K77").Value = 54359375
$ws.Range("L77").Value = 4200815
$ws.Range("M77").Value = -54355007

$ws.Range("H132").Value = 2677.8604
$ws.Range("I132").Value = 1890.7969
$ws.Range("J132").Value = 4967.5
$ws.Range("K132").Value = 5672.3907
$ws.Range("L132").Value = 14902.5
$ws.Range("M132").Value = -3142.3907

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1803
$ws.Range("I86").Value = 1785.091
$ws.Range("J86").Value = 2000
$ws.Range("K86").Value = 1785.091
$ws.Range("L86").Value = 2000
$ws.Range("M86").Value = -662.0909999999999

$ws.Range("H89").Value = 1803
$ws.Range("I89").Value = 1785.091
$ws.Range("J89").Value = 2000
$ws.Range("K89").Value = 8925.455
$ws.Range("L89").Value = 10000
$ws.Range("M89").Value = -3309.455

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 2613.9167
$ws.Range("I58").Value = 1422.1578
$ws.Range("J58").Value = 7142.6
$ws.Range("K58").Value = 1422.1578
$ws.Range("L58").Value = 7142.6
$ws.Range("M58").Value = -1219.1578

$ws.Range("H107").Value = 1816.84
$ws.Range("I107").Value = 843.1667
$ws.Range("J107").Value = 2715.6155
$ws.Range("K107").Value = 843.1667
$ws.Range("L107").Value = 2715.6155
$ws.Range("M107").Value = 1076.8333

$ws.Range("H134").Value = 3772.7693
$ws.Range("I134").Value = 2104.75
$ws.Range("J134").Value = 9332.833000000001
$ws.Range("K134").Value = 6314.25
$ws.Range("L134").Value = 27998.499
$ws.Range("M134").Value = -3779.25

$ws.Range("H135").Value = 94999.664
$ws.Range("I135").Value = 70000
$ws.Range("J135").Value = 107499.5
$ws.Range("K135").Value = 70000
$ws.Range("L135").Value = 107499.5
$ws.Range("M135").Value = -64930
$ws.Range("N135").Value = -117639.5

$ws.Range("H136").Value = 2613.9167
$ws.Range("I136").Value = 1422.1578
$ws.Range("J136").Value = 7142.6
$ws.Range("K136").Value = 4266.4734
$ws.Range("L136").Value = 21427.8
$ws.Range("M136").Value = -1716.4734

$ws.Range("H141").Value = 341332.28
$ws.Range("I141").Value = 99000
$ws.Range("J141").Value = 381721
$ws.Range("K141").Value = 99000
$ws.Range("L141").Value = 381721
$ws.Range("M141").Value = -93820
$ws.Range("N141").Value = -392081

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H15").Value = 236
$ws.Range("I15").Value = 240.83333
$ws.Range("J15").Value = 234.79167
$ws.Range("K15").Value = 722.49999
$ws.Range("L15").Value = 704.37501
$ws.Range("M15").Value = -582.49999

$ws.Range("H24").Value = 3763.8
$ws.Range("I24").Value = 4304.75
$ws.Range("J24").Value = 1600
$ws.Range("K24").Value = 12914.25
$ws.Range("L24").Value = 4800
$ws.Range("M24").Value = -12684.25

$ws.Range("H25").Value = 103.52941
$ws.Range("I25").Value = 109.21429
$ws.Range("J25").Value = 77
$ws.Range("K25").Value = 327.64287
$ws.Range("L25").Value = 231
$ws.Range("M25").Value = -158.64287
$ws.Range("N25").Value = -569

$ws.Range("H30").Value = 103.52941
$ws.Range("I30").Value = 109.21429
$ws.Range("J30").Value = 77
$ws.Range("K30").Value = 327.64287
$ws.Range("L30").Value = 231
$ws.Range("M30").Value = -225.64287
$ws.Range("N30").Value = -435

$ws.Range("H34").Value = 9333.333000000001
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 9333.333000000001
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 27999.999
$ws.Range("N34").Value = -28167.999

$ws.Range("H50").Value = 428.57144
$ws.Range("I50").Value = 299.8
$ws.Range("J50").Value = 456.56522
$ws.Range("K50").Value = 899.4000000000001
$ws.Range("L50").Value = 1369.69566
$ws.Range("M50").Value = -418.4000000000001

$ws.Range("H53").Value = 428.57144
$ws.Range("I53").Value = 299.8
$ws.Range("J53").Value = 456.56522
$ws.Range("K53").Value = 899.4000000000001
$ws.Range("L53").Value = 1369.69566
$ws.Range("M53").Value = -418.4000000000001

$ws.Range("H80").Value = 4278.5713
$ws.Range("I80").Value = 4402
$ws.Range("J80").Value = 4269.077
$ws.Range("K80").Value = 13206
$ws.Range("L80").Value = 12807.231
$ws.Range("M80").Value = -12270
$ws.Range("N80").Value = -14679.231

$ws.Range("H83").Value = 4278.5713
$ws.Range("I83").Value = 4402
$ws.Range("J83").Value = 4269.077
$ws.Range("K83").Value = 39618
$ws.Range("L83").Value = 38421.693
$ws.Range("M83").Value = -34938
$ws.Range("N83").Value = -47781.693

$ws.Range("H107").Value = 596.89655
$ws.Range("I107").Value = 474
$ws.Range("J107").Value = 771
$ws.Range("K107").Value = 1422
$ws.Range("L107").Value = 2313
$ws.Range("M107").Value = 498
$ws.Range("N107").Value = -6153

$ws.Range("H131").Value = 5065.2104
$ws.Range("I131").Value = 4485.231
$ws.Range("J131").Value = 5236.5684
$ws.Range("K131").Value = 13455.693
$ws.Range("L131").Value = 15709.7052
$ws.Range("M131").Value = -8415.692999999999
$ws.Range("N131").Value = -25789.7052

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H69").Value = 112000
$ws.Range("I69").Value = 0
$ws.Range("J69").Value = 112000
$ws.Range("K69").Value = 0
$ws.Range("L69").Value = 112000
$ws.Range("N69").Value = -113498

$ws.Range("H70").Value = 6654.3335
$ws.Range("I70").Value = 5497.5
$ws.Range("J70").Value = 7579.8
$ws.Range("K70").Value = 5497.5
$ws.Range("L70").Value = 7579.8
$ws.Range("M70").Value = -5227.5

$ws.Range("H72").Value = 112000
$ws.Range("I72").Value = 0
$ws.Range("J72").Value = 112000
$ws.Range("K72").Value = 0
$ws.Range("L72").Value = 336000
$ws.Range("N72").Value = -343488

$ws.Range("H73").Value = 6654.3335
$ws.Range("I73").Value = 5497.5
$ws.Range("J73").Value = 7579.8
$ws.Range("K73").Value = 5497.5
$ws.Range("L73").Value = 7579.8
$ws.Range("M73").Value = -4561.5

$ws.Range("H102").Value = 1921.0834
$ws.Range("I102").Value = 1508.0385
$ws.Range("J102").Value = 2995
$ws.Range("K102").Value = 1508.0385
$ws.Range("L102").Value = 2995
$ws.Range("M102").Value = 113.9614999999999

$ws.Range("H132").Value = 25006034
$ws.Range("I132").Value = 30306454
$ws.Range("J132").Value = 18335.857
$ws.Range("K132").Value = 90919362
$ws.Range("L132").Value = 55007.571
$ws.Range("M132").Value = -90916832

$ws.Range("H133").Value = 89245
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 89245
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 89245
$ws.Range("N133").Value = -99365

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 562
$ws.Range("I16").Value = 562
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 562
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -392

$ws.Range("H22").Value = 1467.3334
$ws.Range("I22").Value = 1267.3334
$ws.Range("J22").Value = 1667.3334
$ws.Range("K22").Value = 1267.3334
$ws.Range("L22").Value = 1667.3334
$ws.Range("M22").Value = -972.3334
$ws.Range("N22").Value = -2257.3334

$ws.Range("H27").Value = 1467.3334
$ws.Range("I27").Value = 1267.3334
$ws.Range("J27").Value = 1667.3334
$ws.Range("K27").Value = 1267.3334
$ws.Range("L27").Value = 1667.3334
$ws.Range("M27").Value = -1160.3334
$ws.Range("N27").Value = -1881.3334

$ws.Range("H100").Value = 3720.182
$ws.Range("I100").Value = 2730.75
$ws.Range("J100").Value = 4285.5713
$ws.Range("K100").Value = 2730.75
$ws.Range("L100").Value = 4285.5713
$ws.Range("M100").Value = -2189.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 68500
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 68500
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 68500
$ws.Range("N70").Value = -69130

$ws.Range("H73").Value = 68500
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 68500
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 68500
$ws.Range("N73").Value = -70684

$ws.Range("H81").Value = 17000
$ws.Range("I81").Value = 4000.5
$ws.Range("J81").Value = 29999.5
$ws.Range("K81").Value = 8001
$ws.Range("L81").Value = 59999
$ws.Range("M81").Value = -6940
$ws.Range("N81").Value = -62121

$ws.Range("H84").Value = 17000
$ws.Range("I84").Value = 4000.5
$ws.Range("J84").Value = 29999.5
$ws.Range("K84").Value = 40005
$ws.Range("L84").Value = 299995
$ws.Range("M84").Value = -34701
$ws.Range("N84").Value = -310603

$ws.Range("H95").Value = 87125
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 87125
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 87125
$ws.Range("N95").Value = -92617

$ws.Range("H132").Value = 2123.2307
$ws.Range("I132").Value = 1440.5
$ws.Range("J132").Value = 4399
$ws.Range("K132").Value = 4321.5
$ws.Range("L132").Value = 13197
$ws.Range("M132").Value = -1791.5

$ws.Range("H136").Value = 4392.2856
$ws.Range("I136").Value = 4429.9443
$ws.Range("J136").Value = 4166.3335
$ws.Range("K136").Value = 13289.8329
$ws.Range("L136").Value = 12499.0005
$ws.Range("M136").Value = -10739.8329

$ws.Range("H138").Value = 50143
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 50143
$ws.Range("K138").Value = 0
$ws.Range("L138").Value = 50143
$ws.Range("N138").Value = -60423
